$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.719288
$ws.Range("H2").Value = 11.157864
$ws.Range("I2").Value = 0.04235839908674209
$ws.Range("J2").Value = 0.04235839908674209
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3389413333333333
$ws.Range("N2").Value = 1.016824
$ws.Range("Q2").Value = 1.260620433770667
$ws.Range("R2").Value = 11.345583903936
$ws.Range("S2").Value = 0.04235839908674209
$ws.Range("T2").Value = 0.04235839908674209

# Row 3
$ws.Range("I3").Value = 0.2979256989470644
$ws.Range("J3").Value = 0.2979256989470644
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3389413333333333
$ws.Range("N3").Value = 1.016824
$ws.Range("Q3").Value = 8.86651129257689
$ws.Range("R3").Value = 79.79860163319201
$ws.Range("S3").Value = 0.2979256989470644
$ws.Range("T3").Value = 0.2979256989470644

# Row 4
$ws.Range("G4").Value = 6.299630666666666
$ws.Range("H4").Value = 18.898892
$ws.Range("I4").Value = 0.07174552491706633
$ws.Range("J4").Value = 0.07174552491706633
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3389413333333333
$ws.Range("N4").Value = 1.016824
$ws.Range("Q4").Value = 2.135205217667556
$ws.Range("R4").Value = 19.216846959008
$ws.Range("S4").Value = 0.07174552491706633
$ws.Range("T4").Value = 0.07174552491706633

# Row 5
$ws.Range("G5").Value = 51.62686066666667
$ws.Range("H5").Value = 154.880582
$ws.Range("I5").Value = 0.5879703770491272
$ws.Range("J5").Value = 0.5879703770491272
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3389413333333333
$ws.Range("N5").Value = 1.016824
$ws.Range("Q5").Value = 17.49847699017422
$ws.Range("R5").Value = 157.486292911568
$ws.Range("S5").Value = 0.5879703770491272
$ws.Range("T5").Value = 0.5879703770491272
